# Updates the "Avverkningsanmälningar" sheet:
#  1. Refresh the "Förändrad" (column C) date for every existing data row
#     (rows 2-485) from 2023-09-23 (45192) to 2023-10-03 (45202).
#  2. Append three new report rows (486-488) pulled in by the data refresh.
#  3. Make sure row 485 carries an explicit row height (as the newly
#     appended rows do), matching the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD"

# ---------------------------------------------------------------------
# 1. Bulk-update the "Förändrad" column for all existing data rows.
# ---------------------------------------------------------------------
$ws.Range("C2:C485").Value2 = 45202

# ---------------------------------------------------------------------
# 2. Give row 485 an explicit height, same as the newly appended rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(485).RowHeight = 15

# ---------------------------------------------------------------------
# 3. Append the three new rows of data.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 486; A = "A 46528-2023"; B = 45197; C = 45202; D = "JÖNKÖPINGS LÄN"; E = "EKSJÖ"; G = 0.7 },
    @{ Row = 487; A = "A 46309-2023"; B = 45197; C = 45202; D = "JÖNKÖPINGS LÄN"; E = "EKSJÖ"; G = 1   },
    @{ Row = 488; A = "A 47022-2023"; B = 45201; C = 45202; D = "JÖNKÖPINGS LÄN"; E = "EKSJÖ"; G = 3.7 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # The sheet stamps an explicit row height on every row except the
    # very last one in the data range, so mirror that here too.
    if ($rowIndex -ne 488) {
        $ws.Rows.Item($rowIndex).RowHeight = 15
    }

    $ws.Cells.Item($rowIndex, 1).Value2 = $r.A

    $ws.Cells.Item($rowIndex, 2).Value2 = $r.B
    $ws.Cells.Item($rowIndex, 2).NumberFormat = $dateFormat

    $ws.Cells.Item($rowIndex, 3).Value2 = $r.C
    $ws.Cells.Item($rowIndex, 3).NumberFormat = $dateFormat

    $ws.Cells.Item($rowIndex, 4).Value2 = $r.D
    $ws.Cells.Item($rowIndex, 5).Value2 = $r.E

    # F (Markägare) intentionally left blank, matching the source data.

    $ws.Cells.Item($rowIndex, 7).Value2 = $r.G

    # H..Q (columns 8-17) are all zero counts.
    for ($c = 8; $c -le 17; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value2 = 0
    }

    # R (Artnamn) stays empty but keeps the wrap-text styling used
    # throughout the rest of the sheet.
    $ws.Cells.Item($rowIndex, 18).Value2 = ""
    $ws.Cells.Item($rowIndex, 18).WrapText = $true
}
